$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# --- Row 2: replace the single-model (CESM2-CAM6) summary row with the
#     ensemble "Model Output" row. Columns B:G keep their existing style
#     (s=5 for data cells, s=8 for the label) and only the values/text
#     change; columns H:N are no longer part of this row and are cleared.
$ws.Range("A2").Value = "Model Output"
$ws.Range("B2").Value = 50.17197
$ws.Range("C2").Value = 247.09817
$ws.Range("D2").Value = 262.99485
$ws.Range("E2").Value = 321.183526
$ws.Range("F2").Value = 98.55575
$ws.Range("G2").Value = 239.06726
$ws.Range("H2:N2").ClearContents()

# --- Row 6: previously blank (style-only) spacer row now holds the
#     MRI-ESM2 model row, mirroring the old row 2 layout/values.
$ws.Range("A6").Value = "MRI-ESM2"
$ws.Range("B6").Value = 4.036277054379609
$ws.Range("C6").Value = 342.7713118934579
$ws.Range("D6").Value = 98.78611701806761
$ws.Range("E6").Value = 53.97546924823823
$ws.Range("F6").Value = 239.9489180602529
$ws.Range("G6").Value = 264.3077785917636
$ws.Range("H6").Value = 192.4933453275627
$ws.Range("I6").Value = 245.8657797562242
$ws.Range("J6").Value = 26.85690608577296
$ws.Range("K6").Value = 32.62033959045065
$ws.Range("L6").Value = 403.4433355451918
$ws.Range("M6").Value = 347.1772706961048
$ws.Range("N6").Value = 319.8823832549311
